$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All text-like values are prefixed with a leading apostrophe so Excel
# keeps them as Text (matching the original inlineStr cell type) instead
# of auto-converting numeric-looking strings (e.g. "245.56") into numbers
# or stripping formatting like trailing zeros / leading zeros in exponents.

# Row 2
$ws.Range("D2").Value = "'30.541.67"
$ws.Range("E2").Value = "'  +0.17%  "

# Row 3
$ws.Range("D3").Value = "'1.919.08"
$ws.Range("E3").Value = "'  -0.23%  "

# Row 4
$ws.Range("E4").Value = "'  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'245.56"
$ws.Range("E5").Value = "'  +1.03%  "

# Row 6
$ws.Range("E6").Value = "'  +0.00%  "

# Row 7
$ws.Range("D7").Value = "'0.4791"
$ws.Range("E7").Value = "'  +1.71%  "

# Row 8
$ws.Range("D8").Value = "'0.2884"
$ws.Range("E8").Value = "'  +0.17%  "

# Row 9
$ws.Range("D9").Value = "'0.06721"
$ws.Range("E9").Value = "'  -0.45%  "

# Row 10
$ws.Range("D10").Value = "'109.82"
$ws.Range("E10").Value = "'  +3.34%  "

# Row 11
$ws.Range("D11").Value = "'19.15"
$ws.Range("E11").Value = "'  +4.66%  "

# Row 12
$ws.Range("D12").Value = "'1.917.05"
$ws.Range("E12").Value = "'  -0.29%  "

# Row 13
$ws.Range("D13").Value = "'0.07570"
$ws.Range("E13").Value = "'  -2.31%  "

# Row 14
$ws.Range("D14").Value = "'5.262"
$ws.Range("E14").Value = "'  -0.79%  "

# Row 15
$ws.Range("D15").Value = "'0.6680"
$ws.Range("E15").Value = "'  +1.36%  "

# Row 16
$ws.Range("D16").Value = "'298.97"
$ws.Range("E16").Value = "'  +2.57%  "

# Row 17
$ws.Range("D17").Value = "'30.534.77"
$ws.Range("E17").Value = "'  +0.13%  "

# Row 18
$ws.Range("D18").Value = "'13.04"
$ws.Range("E18").Value = "'  +0.82%  "

# Row 19
$ws.Range("D19").Value = "'5.594"
$ws.Range("E19").Value = "'  +6.27%  "

# Row 20
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "'  +0.02%  "

# Row 21
$ws.Range("D21").Value = "'0.000007569"
$ws.Range("E21").Value = "'  -0.23%  "

# Row 22
$ws.Range("D22").Value = "'2.165.20"
$ws.Range("E22").Value = "'  +0.57%  "

# Row 23
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "'  +0.17%  "

# Row 24
$ws.Range("D24").Value = "'6.428"
$ws.Range("E24").Value = "'  +3.58%  "

# Row 25
$ws.Range("D25").Value = "'9.483"
$ws.Range("E25").Value = "'  +1.23%  "

# Row 26
$ws.Range("D26").Value = "'164.61"
$ws.Range("E26").Value = "'  -2.58%  "

# Row 27
$ws.Range("D27").Value = "'20.29"
$ws.Range("E27").Value = "'  -5.21%  "

# Row 28
$ws.Range("D28").Value = "'2.115"
$ws.Range("E28").Value = "'  +0.26%  "

# Row 29
$ws.Range("D29").Value = "'0.1077"
$ws.Range("E29").Value = "'  +0.80%  "

# Row 30
$ws.Range("D30").Value = "'1.395"
$ws.Range("E30").Value = "'  +2.02%  "

# Row 31
$ws.Range("E31").Value = "'  -0.34%  "

# Row 32
$ws.Range("D32").Value = "'4.037"
$ws.Range("E32").Value = "'  +1.28%  "

# Row 33
$ws.Range("D33").Value = "'0.04998"
$ws.Range("E33").Value = "'  -0.57%  "

# Row 34
$ws.Range("D34").Value = "'0.7366"
$ws.Range("E34").Value = "'  -0.69%  "

# Row 35
$ws.Range("D35").Value = "'1.138"
$ws.Range("E35").Value = "'  -1.38%  "

# Row 36
$ws.Range("D36").Value = "'0.9996"
$ws.Range("E36").Value = "'  +0.00%  "

# Row 37
$ws.Range("B37").Value = "'HuobiToken"
$ws.Range("C37").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.724"
$ws.Range("E37").Value = "'  -0.13%  "

# Row 38
$ws.Range("B38").Value = "'VeChain"
$ws.Range("C38").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02036"
$ws.Range("E38").Value = "'  -3.16%  "

# Row 39
$ws.Range("D39").Value = "'2.686"
$ws.Range("E39").Value = "'  +0.08%  "

# Row 40
$ws.Range("D40").Value = "'110.78"
$ws.Range("E40").Value = "'  +0.52%  "

# Row 41
$ws.Range("D41").Value = "'2.019"
$ws.Range("E41").Value = "'  -2.42%  "

# Row 42
$ws.Range("D42").Value = "'0.4429"
$ws.Range("E42").Value = "'  +4.05%  "

# Row 43
$ws.Range("B43").Value = "'Aave"
$ws.Range("C43").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'72.83"
$ws.Range("E43").Value = "'  +8.08%  "

# Row 44
$ws.Range("B44").Value = "'TrustWalletToken"
$ws.Range("C44").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8616"
$ws.Range("E44").Value = "'  -1.26%  "

# Row 45
$ws.Range("D45").Value = "'5.909"
$ws.Range("E45").Value = "'  +0.99%  "

# Row 46
$ws.Range("E46").Value = "'  +0.03%  "

# Row 47
$ws.Range("D47").Value = "'49.43"
$ws.Range("E47").Value = "'  +0.90%  "

# Row 48
$ws.Range("D48").Value = "'7.269"
$ws.Range("E48").Value = "'  +1.07%  "

# Row 49
$ws.Range("D49").Value = "'9.308"
$ws.Range("E49").Value = "'  +0.34%  "

# Row 50
$ws.Range("D50").Value = "'0.1231"
$ws.Range("E50").Value = "'  +1.25%  "

# Row 51
$ws.Range("E51").Value = "'  +2.74%  "
